$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: update existing values ---
$ws.Range("A2").Value = "FV/001"
$ws.Range("B2").Value = 45843
$ws.Range("C2").Value = "Materiały MO"
$ws.Range("D2").Value = 60000
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "73"
$ws.Range("E2").ClearFormats()

# --- Row 3: new row ---
$ws.Range("A3").Value = "FV/002"
$ws.Range("B3").NumberFormat = "YYYY-MM-DD"
$ws.Range("B3").Value = 45850
$ws.Range("C3").Value = "Energia MO"
$ws.Range("D3").Value = 40000
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "73"
$ws.Range("E3").ClearFormats()

# --- Row 4: new row ---
$ws.Range("A4").Value = "FV/003"
$ws.Range("B4").NumberFormat = "YYYY-MM-DD"
$ws.Range("B4").Value = 45857
$ws.Range("C4").Value = "Lakier LAK"
$ws.Range("D4").Value = 50000
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "10"
$ws.Range("E4").ClearFormats()

# --- Row 5: new row ---
$ws.Range("A5").Value = "FV/004"
$ws.Range("B5").NumberFormat = "YYYY-MM-DD"
$ws.Range("B5").Value = 45863
$ws.Range("C5").Value = "Korekta MO"
$ws.Range("D5").Value = -5000
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "73"
$ws.Range("E5").ClearFormats()
